$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 12 ("Sharing") - TextBox 1: grow the box and fix three typos
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shp = $s12.Shapes.Item(2)

# Grow the text box height (13865629 x 5170646 EMU -> 13865629 x 5632311 EMU)
$shp.Height = 443.4891

$tr = $s12.Shapes.Item(2).TextFrame.TextRange

# 1) "... and download the latest Python 3" -> "... to download and install the latest Python 3"
$found1 = $tr.Find(" and download the latest Python 3")
$found1.Text = " to download and install the latest Python 3"

# 2) "3 - Teach your friends who to run " -> "3 - Teach your friends how to run " (keep curly quote intact)
$full = $tr.Text
$needle = "3 " + [string][char]0x2013 + " Teach your friends who to run"
$idx = $full.IndexOf($needle)
$run2 = $tr.Characters($idx + 1, 35)
$run2.Text = "3 " + [string][char]0x2013 + " Teach your friends how to run " + [string][char]0x201C

# 3) "... install VS Code so they can change your game!" -> "... install Visual Studio Code so they can change your game!"
$found3 = $tr.Find("You can also show them how to install VS Code so they can change your game!")
$found3.Text = "You can also show them how to install Visual Studio Code so they can change your game!"

# ---------------------------------------------------------------------------
# Slide 8 ("A Quick Example") - wormy_mod picture: nudge vertical position
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$pic = $s8.Shapes.Item(3)
$pic.Top = 165.112
